$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the starting totals used in the "Best" block (C column)
$ws.Range("C3").Value = 995739

# Update the starting totals used in the "Local" block (F column)
$ws.Range("F3").Value = 954859

# Update the selected cell to match the saved selection state
$ws.Range("G10").Select()
